$d = $word.ActiveDocument

# --- Edit 1: append the "(This is a change - Version for branch alternate)" note
#     to the first paragraph, in dark-red (C00000) text, split across three runs
#     (matching how the text was actually typed/pasted in stages).
$p1 = $d.Paragraphs(1).Range
$insertPos = $p1.End - 1                      # just before the paragraph mark
$cursor = $d.Range($insertPos, $insertPos)

# two plain trailing spaces, no special formatting
$cursor.InsertAfter("  ")

# first colored run: "(This is a change " + en dash + " Ve"
$run1 = $d.Range($cursor.End, $cursor.End)
$run1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1.Font.Color = 192

# second colored run: "rsion for branch alternate"
$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 192

# third colored run: ")"
$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 192

# --- Edit 2: append a new, otherwise-empty paragraph shaded F9F9F9 after the
#     final "Free at last" paragraph.
$endPos = $d.Content.End
$tail = $d.Range($endPos, $endPos)
$tail.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"F9F9F9`"/></w:pPr></w:p>")

# --- Edit 3: the resave also pruned the unused heading/web-paste styles that
#     nothing in the body references (document has no headings, no pasted-html
#     classes) - remove them too, highest index first so the collection never
#     has to be re-queried by name after a shift.
$unusedStyleIdx = @(18, 17, 16, 15, 14, 13, 12, 11, 10, 3, 2)
foreach ($idx in $unusedStyleIdx) {
    $d.Styles.Item($idx).Delete()
}
